$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F272").Value = 30812
$ws.Range("F273").Value = 31724
$ws.Range("F274").Value = 28111
$ws.Range("G274").Value = 1279
$ws.Range("F275").Value = 30348
$ws.Range("F277").Value = 3396
$ws.Range("F278").Value = 30547
$ws.Range("F279").Value = 42767
$ws.Range("F280").Value = 34793
$ws.Range("G280").Value = 2311
$ws.Range("F281").Value = 46079
$ws.Range("F286").Value = 55194
$ws.Range("F287").Value = 58880
$ws.Range("G287").Value = 3719
$ws.Range("F288").Value = 59280
$ws.Range("F289").Value = 62991
$ws.Range("F290").Value = 17585
$ws.Range("F292").Value = 82470
$ws.Range("F294").Value = 93964
$ws.Range("G294").Value = 4949
$ws.Range("F300").Value = 72572
$ws.Range("F301").Value = 72198
$ws.Range("G301").Value = 5684
$ws.Range("F302").Value = 78626
$ws.Range("G302").Value = 5657
$ws.Range("F307").Value = 75895
$ws.Range("G307").Value = 6399
$ws.Range("F308").Value = 15471
$ws.Range("F309").Value = 77920
$ws.Range("F314").Value = 64368
$ws.Range("G314").Value = 3149
$ws.Range("F315").Value = 56369
$ws.Range("G315").Value = 2628
$ws.Range("F316").Value = 50750
$ws.Range("F317").Value = 63739
$ws.Range("F321").Value = 89350
$ws.Range("G321").Value = 2655
$ws.Range("F322").Value = 109661
$ws.Range("F323").Value = 216839
$ws.Range("F324").Value = 241037
$ws.Range("F325").Value = 766046
$ws.Range("G325").Value = 6456
$ws.Range("F326").Value = 419556
$ws.Range("F327").Value = 224977
$ws.Range("G327").Value = 2725
$ws.Range("F328").Value = 180851
$ws.Range("G328").Value = 2671
$ws.Range("F329").Value = 82999
$ws.Range("G329").Value = 1759
$ws.Range("F330").Value = 72552
$ws.Range("G330").Value = 2084
$ws.Range("F331").Value = 155001
$ws.Range("F332").Value = 457208
$ws.Range("F333").Value = 271711
$ws.Range("G333").Value = 2945
$ws.Range("F334").Value = 196778
$ws.Range("F335").Value = 130955
$ws.Range("G335").Value = 3001
$ws.Range("F336").Value = 102629
$ws.Range("G336").Value = 3344
$ws.Range("F337").Value = 103566
$ws.Range("G337").Value = 2891
$ws.Range("F338").Value = 228389
$ws.Range("F339").Value = 661293
$ws.Range("F341").Value = 291779
$ws.Range("F342").Value = 178740
$ws.Range("G342").Value = 3039
$ws.Range("F343").Value = 133289
$ws.Range("F344").Value = 135483
$ws.Range("F347").Value = 343715
$ws.Range("F349").Value = 159364
$ws.Range("F350").Value = 127094
$ws.Range("F351").Value = 150927
$ws.Range("G351").Value = 2833
$ws.Range("F352").Value = 307409
$ws.Range("G352").Value = 3541
$ws.Range("F353").Value = 725597
$ws.Range("F355").Value = 222047
$ws.Range("F356").Value = 160046
$ws.Range("G356").Value = 2877
$ws.Range("F357").Value = 138506
$ws.Range("G357").Value = 3028
$ws.Range("F358").Value = 158741
$ws.Range("F359").Value = 321219
$ws.Range("F360").Value = 751797
$ws.Range("F362").Value = 229127
$ws.Range("F364").Value = 168375
$ws.Range("F366").Value = 339834
$ws.Range("F369").Value = 234555
$ws.Range("G369").Value = 2601
$ws.Range("F373").Value = 349537
$ws.Range("G373").Value = 2371
$ws.Range("F374").Value = 773231
$ws.Range("F376").Value = 220929
$ws.Range("F378").Value = 157549
$ws.Range("F382").Value = 357995
$ws.Range("F384").Value = 171991
$ws.Range("F387").Value = 351635
$ws.Range("F388").Value = 729298
$ws.Range("F391").Value = 176842
$ws.Range("F393").Value = 306870
$ws.Range("G393").Value = 1227
$ws.Range("F394").Value = 166155
$ws.Range("G394").Value = 633
$ws.Range("F395").Value = 749209
$ws.Range("G395").Value = 1945
$ws.Range("F398").Value = 297915
$ws.Range("G398").Value = 1467
$ws.Range("F399").Value = 200374
$ws.Range("G399").Value = 969
$ws.Range("F400").Value = 149470
$ws.Range("G400").Value = 756
$ws.Range("F401").Value = 272886
$ws.Range("G401").Value = 934
$ws.Range("F402").Value = 712886
$ws.Range("G402").Value = 1366
$ws.Range("F403").Value = 349452
$ws.Range("F404").Value = 224545
$ws.Range("G404").Value = 905
$ws.Range("F405").Value = 173282
$ws.Range("G405").Value = 692
$ws.Range("F406").Value = 170091
$ws.Range("G406").Value = 673
$ws.Range("F407").Value = 155874
$ws.Range("G407").Value = 666

# New rows appended at the end
$ws.Cells.Item(408,1).Value = 44302
$ws.Cells.Item(408,2).Value = 375336
$ws.Cells.Item(408,3).Value = 13377
$ws.Cells.Item(408,4).Value = 750
$ws.Cells.Item(408,5).Value = 11043
$ws.Cells.Item(408,6).Value = 284565
$ws.Cells.Item(408,7).Value = 1087

$ws.Cells.Item(409,1).Value = 44303
$ws.Cells.Item(409,2).Value = 375974
$ws.Cells.Item(409,3).Value = 17057
$ws.Cells.Item(409,4).Value = 638
$ws.Cells.Item(409,5).Value = 11106
$ws.Cells.Item(409,6).Value = 587476
$ws.Cells.Item(409,7).Value = 2028

"Update complete"